$wb = $excel.ActiveWorkbook

# --- Text content updates -------------------------------------------------

# "data" sheet: rename title "xyz" -> "Xyz"
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("A1").Value = "Xyz"

# "day 1" sheet: "Day 1: Chicken ..." -> "Day 1: Legs ..."
$wsDay1 = $wb.Worksheets.Item("day 1")
$wsDay1.Range("A1").Value = "Day 1: Legs ĄĆĘÓŁŃŻŹąćęółńżź ĄĆĘÓŁŃŻŹąćęółńżź"

# "day 3" sheet: "Day 3: popo" -> "Day 3: ABS and arms"
$wsDay3 = $wb.Worksheets.Item("day 3")
$wsDay3.Range("A1").Value = "Day 3: ABS and arms"

# --- Column width on "day 2" sheet -----------------------------------------

$wsDay2 = $wb.Worksheets.Item("day 2")
$wsDay2.Columns.Item(1).ColumnWidth = 28.94

# --- Reset stray cell selections on sheets whose cursor moved back to A1 ---
# (data / day 1 / day 3 no longer keep their old mid-sheet selection)

$excel.Goto($wsDay1.Range("A1"), $false)
$excel.Goto($wsDay3.Range("A1"), $false)

# --- Active sheet / tab selection ------------------------------------------
# Previously "day 3" was the active tab; now "data" should be active.

$wsData.Activate()
$wsData.Range("A1").Select()
